# Applies the FlashScore 2024-12-13 odds update described in the commit diff.
# All changes are plain numeric overwrites of existing "Odd_*" cells in Sheet1
# (rows 2,3,4,6,7,8,9,10,12,15,16,17,20,44), leaving every other cell untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 1.44
$ws.Cells.Item(2, 8).Value = 4.5
$ws.Cells.Item(2, 9).Value = 7.5
$ws.Cells.Item(2, 10).Value = 1.95
$ws.Cells.Item(2, 11).Value = 2.38
$ws.Cells.Item(2, 12).Value = 6.5
$ws.Cells.Item(2, 15).Value = 1.25
$ws.Cells.Item(2, 16).Value = 4
$ws.Cells.Item(2, 17).Value = 1.8
$ws.Cells.Item(2, 18).Value = 2
$ws.Cells.Item(2, 21).Value = 1.95
$ws.Cells.Item(2, 22).Value = 1.8
$ws.Cells.Item(2, 24).Value = 7
$ws.Cells.Item(2, 26).Value = 9.5
$ws.Cells.Item(2, 29).Value = 12
$ws.Cells.Item(2, 30).Value = 8.5
$ws.Cells.Item(2, 31).Value = 19
$ws.Cells.Item(2, 32).Value = 51
$ws.Cells.Item(2, 33).Value = 351
$ws.Cells.Item(2, 34).Value = 17
$ws.Cells.Item(2, 35).Value = 34
$ws.Cells.Item(2, 36).Value = 21
$ws.Cells.Item(2, 37).Value = 81
$ws.Cells.Item(2, 40).Value = 3.4
$ws.Cells.Item(2, 41).Value = 7
$ws.Cells.Item(2, 43).Value = 21
$ws.Cells.Item(2, 45).Value = 126
$ws.Cells.Item(2, 47).Value = 9
$ws.Cells.Item(2, 48).Value = 51
$ws.Cells.Item(2, 49).Value = 8
$ws.Cells.Item(2, 50).Value = 34
$ws.Cells.Item(2, 54).Value = 301
$ws.Cells.Item(3, 15).Value = 1.25
$ws.Cells.Item(3, 16).Value = 4
$ws.Cells.Item(3, 17).Value = 1.87
$ws.Cells.Item(3, 18).Value = 2.03
$ws.Cells.Item(4, 7).Value = 2.55
$ws.Cells.Item(4, 9).Value = 3.25
$ws.Cells.Item(4, 10).Value = 3.4
$ws.Cells.Item(4, 24).Value = 10
$ws.Cells.Item(4, 27).Value = 26
$ws.Cells.Item(4, 34).Value = 7
$ws.Cells.Item(4, 49).Value = 5
$ws.Cells.Item(4, 50).Value = 21
$ws.Cells.Item(4, 51).Value = 41
$ws.Cells.Item(6, 7).Value = 3.1
$ws.Cells.Item(6, 9).Value = 2.35
$ws.Cells.Item(6, 10).Value = 4
$ws.Cells.Item(6, 12).Value = 3.2
$ws.Cells.Item(6, 15).Value = 1.5
$ws.Cells.Item(6, 16).Value = 2.63
$ws.Cells.Item(6, 17).Value = 2.5
$ws.Cells.Item(6, 18).Value = 1.53
$ws.Cells.Item(6, 23).Value = 8
$ws.Cells.Item(6, 24).Value = 15
$ws.Cells.Item(6, 29).Value = 7.5
$ws.Cells.Item(6, 35).Value = 10
$ws.Cells.Item(6, 38).Value = 21
$ws.Cells.Item(6, 39).Value = 34
$ws.Cells.Item(6, 40).Value = 5
$ws.Cells.Item(6, 45).Value = 251
$ws.Cells.Item(6, 47).Value = 8.5
$ws.Cells.Item(6, 51).Value = 26
$ws.Cells.Item(6, 56).Value = 151
$ws.Cells.Item(7, 7).Value = 6.25
$ws.Cells.Item(7, 8).Value = 3.7
$ws.Cells.Item(7, 9).Value = 1.62
$ws.Cells.Item(7, 10).Value = 7
$ws.Cells.Item(7, 12).Value = 2.3
$ws.Cells.Item(7, 21).Value = 2.5
$ws.Cells.Item(7, 22).Value = 1.5
$ws.Cells.Item(7, 27).Value = 51
$ws.Cells.Item(7, 32).Value = 101
$ws.Cells.Item(7, 40).Value = 7.5
$ws.Cells.Item(7, 44).Value = 251
$ws.Cells.Item(7, 48).Value = 81
$ws.Cells.Item(8, 7).Value = 2.45
$ws.Cells.Item(8, 9).Value = 3.1
$ws.Cells.Item(8, 10).Value = 3.25
$ws.Cells.Item(8, 26).Value = 23
$ws.Cells.Item(8, 37).Value = 34
$ws.Cells.Item(8, 49).Value = 5
$ws.Cells.Item(9, 7).Value = 2.4
$ws.Cells.Item(9, 8).Value = 3
$ws.Cells.Item(9, 9).Value = 3.3
$ws.Cells.Item(9, 11).Value = 1.83
$ws.Cells.Item(9, 13).Value = 1.13
$ws.Cells.Item(9, 14).Value = 6
$ws.Cells.Item(9, 24).Value = 9.5
$ws.Cells.Item(9, 35).Value = 15
$ws.Cells.Item(9, 41).Value = 15
$ws.Cells.Item(9, 53).Value = 151
$ws.Cells.Item(10, 7).Value = 1.95
$ws.Cells.Item(10, 9).Value = 4.5
$ws.Cells.Item(10, 10).Value = 2.88
$ws.Cells.Item(10, 12).Value = 5
$ws.Cells.Item(10, 15).Value = 1.5
$ws.Cells.Item(10, 16).Value = 2.5
$ws.Cells.Item(10, 17).Value = 2.6
$ws.Cells.Item(10, 18).Value = 1.48
$ws.Cells.Item(10, 19).Value = 1.57
$ws.Cells.Item(10, 20).Value = 2.25
$ws.Cells.Item(10, 21).Value = 2.2
$ws.Cells.Item(10, 22).Value = 1.62
$ws.Cells.Item(10, 23).Value = 5.5
$ws.Cells.Item(10, 24).Value = 8
$ws.Cells.Item(10, 31).Value = 19
$ws.Cells.Item(10, 34).Value = 9
$ws.Cells.Item(10, 36).Value = 15
$ws.Cells.Item(10, 41).Value = 12
$ws.Cells.Item(10, 46).Value = 2.25
$ws.Cells.Item(10, 47).Value = 9.5
$ws.Cells.Item(10, 50).Value = 26
$ws.Cells.Item(10, 54).Value = 401
$ws.Cells.Item(12, 7).Value = 2.35
$ws.Cells.Item(12, 9).Value = 2.9
$ws.Cells.Item(12, 10).Value = 3
$ws.Cells.Item(12, 12).Value = 3.6
$ws.Cells.Item(12, 13).Value = 1.05
$ws.Cells.Item(12, 14).Value = 11
$ws.Cells.Item(12, 15).Value = 1.3
$ws.Cells.Item(12, 16).Value = 3.4
$ws.Cells.Item(12, 17).Value = 2.03
$ws.Cells.Item(12, 18).Value = 1.83
$ws.Cells.Item(12, 23).Value = 8
$ws.Cells.Item(12, 39).Value = 34
$ws.Cells.Item(12, 40).Value = 4.33
$ws.Cells.Item(12, 51).Value = 26
$ws.Cells.Item(12, 53).Value = 81
$ws.Cells.Item(15, 7).Value = 3
$ws.Cells.Item(15, 9).Value = 2.3
$ws.Cells.Item(15, 10).Value = 3.5
$ws.Cells.Item(15, 11).Value = 2.2
$ws.Cells.Item(15, 13).Value = 1.05
$ws.Cells.Item(15, 14).Value = 11
$ws.Cells.Item(15, 15).Value = 1.29
$ws.Cells.Item(15, 16).Value = 3.5
$ws.Cells.Item(15, 17).Value = 1.95
$ws.Cells.Item(15, 18).Value = 1.85
$ws.Cells.Item(15, 19).Value = 1.36
$ws.Cells.Item(15, 20).Value = 3
$ws.Cells.Item(15, 24).Value = 15
$ws.Cells.Item(15, 25).Value = 11
$ws.Cells.Item(15, 27).Value = 23
$ws.Cells.Item(15, 31).Value = 13
$ws.Cells.Item(15, 32).Value = 41
$ws.Cells.Item(15, 36).Value = 9.5
$ws.Cells.Item(15, 38).Value = 19
$ws.Cells.Item(15, 42).Value = 23
$ws.Cells.Item(15, 44).Value = 67
$ws.Cells.Item(15, 45).Value = 151
$ws.Cells.Item(15, 46).Value = 3
$ws.Cells.Item(15, 47).Value = 7.5
$ws.Cells.Item(15, 49).Value = 4.5
$ws.Cells.Item(15, 54).Value = 126
$ws.Cells.Item(16, 11).Value = 3.5
$ws.Cells.Item(16, 15).Value = 1.06
$ws.Cells.Item(16, 16).Value = 10
$ws.Cells.Item(16, 21).Value = 1.91
$ws.Cells.Item(16, 22).Value = 1.8
$ws.Cells.Item(16, 23).Value = 13
$ws.Cells.Item(16, 24).Value = 8.5
$ws.Cells.Item(16, 26).Value = 7.5
$ws.Cells.Item(16, 28).Value = 26
$ws.Cells.Item(16, 29).Value = 34
$ws.Cells.Item(16, 31).Value = 29
$ws.Cells.Item(16, 32).Value = 67
$ws.Cells.Item(16, 33).Value = 251
$ws.Cells.Item(16, 37).Value = 201
$ws.Cells.Item(16, 38).Value = 81
$ws.Cells.Item(16, 40).Value = 3.6
$ws.Cells.Item(16, 44).Value = 26
$ws.Cells.Item(16, 45).Value = 81
$ws.Cells.Item(16, 51).Value = 41
$ws.Cells.Item(16, 52).Value = 301
$ws.Cells.Item(16, 53).Value = 201
$ws.Cells.Item(16, 54).Value = 301
$ws.Cells.Item(17, 15).Value = 1.13
$ws.Cells.Item(17, 16).Value = 6
$ws.Cells.Item(17, 17).Value = 1.44
$ws.Cells.Item(17, 18).Value = 2.7
$ws.Cells.Item(20, 10).Value = 2.75
$ws.Cells.Item(20, 13).Value = 1.08
$ws.Cells.Item(20, 14).Value = 8
$ws.Cells.Item(20, 15).Value = 1.4
$ws.Cells.Item(20, 16).Value = 3
$ws.Cells.Item(20, 19).Value = 1.5
$ws.Cells.Item(20, 20).Value = 2.5
$ws.Cells.Item(20, 21).Value = 2
$ws.Cells.Item(20, 22).Value = 1.73
$ws.Cells.Item(20, 28).Value = 34
$ws.Cells.Item(20, 29).Value = 8
$ws.Cells.Item(20, 30).Value = 6
$ws.Cells.Item(20, 33).Value = 401
$ws.Cells.Item(20, 36).Value = 15
$ws.Cells.Item(20, 46).Value = 2.5
$ws.Cells.Item(20, 53).Value = 126
$ws.Cells.Item(20, 54).Value = 301
$ws.Cells.Item(20, 55).Value = 151
$ws.Cells.Item(44, 17).Value = 2.15
$ws.Cells.Item(44, 18).Value = 1.67

Write-Host "Applied 201 cell updates."
